$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.810.71"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.290.99"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.69"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.56"
$ws.Range("E6").Value = "  -2.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.58"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.68"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "2.648.46"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "2.289.42"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "42.745.27"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -4.43%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.75"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.02"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.02"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.12"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.22"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.85"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.04"
$ws.Range("E36").Value = "  -6.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").Value = "2.009.92"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.07"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.16"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").Value = "2.514.83"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.02"
$ws.Range("E51").Value = "  -3.16%  "
